# Updates the "想去人数" (want-to-go count) figures in column F across the
# 展览 (Exhibition), 演出 (Performance) and 全部类型 (All types) sheets to match
# the refreshed data pulled at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 333
$ws1.Range("F3").Value = 280
$ws1.Range("F4").Value = 1227
$ws1.Range("F5").Value = 275
$ws1.Range("F8").Value = 42
$ws1.Range("F9").Value = 137
$ws1.Range("F10").Value = 3428
$ws1.Range("F11").Value = 122
$ws1.Range("F15").Value = 53
$ws1.Range("F16").Value = 587
$ws1.Range("F18").Value = 710
$ws1.Range("F20").Value = 116
$ws1.Range("F22").Value = 55
$ws1.Range("F23").Value = 63
$ws1.Range("F24").Value = 2542
$ws1.Range("F25").Value = 5050
$ws1.Range("F31").Value = 2222
$ws1.Range("F33").Value = 488
$ws1.Range("F35").Value = 101
$ws1.Range("F36").Value = 168
$ws1.Range("F39").Value = 786
$ws1.Range("F42").Value = 33

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 69

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 333
$ws4.Range("F3").Value = 280
$ws4.Range("F4").Value = 1227
$ws4.Range("F5").Value = 275
$ws4.Range("F8").Value = 42
$ws4.Range("F9").Value = 137
$ws4.Range("F10").Value = 3428
$ws4.Range("F11").Value = 122
$ws4.Range("F14").Value = 69
$ws4.Range("F16").Value = 53
$ws4.Range("F17").Value = 587
$ws4.Range("F19").Value = 710
$ws4.Range("F21").Value = 116
$ws4.Range("F23").Value = 55
$ws4.Range("F24").Value = 63
$ws4.Range("F25").Value = 2542
$ws4.Range("F26").Value = 5050
$ws4.Range("F32").Value = 2222
$ws4.Range("F34").Value = 488
$ws4.Range("F36").Value = 101
$ws4.Range("F37").Value = 168
$ws4.Range("F40").Value = 786
$ws4.Range("F43").Value = 33
